$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.691.02"
$ws.Range("E2").Value = "  +0.86%  "

$ws.Range("D3").Value = "1.645.90"
$ws.Range("E3").Value = "  +1.15%  "

$ws.Range("E4").Value = "  +0.32%  "

$ws.Range("D5").Value = "'215.78"
$ws.Range("E5").Value = "  +1.21%  "

$ws.Range("E6").Value = "  +1.49%  "

$ws.Range("E7").Value = "  +0.27%  "

$ws.Range("E8").Value = "  +1.20%  "

$ws.Range("D9").Value = "'0.0626"
$ws.Range("E9").Value = "  +0.45%  "

$ws.Range("D10").Value = "'19.11"
$ws.Range("E10").Value = "  +0.89%  "

$ws.Range("D11").Value = "'0.0846"
$ws.Range("E11").Value = "  +0.10%  "

$ws.Range("D12").Value = "1.876.47"
$ws.Range("E12").Value = "  +1.22%  "

$ws.Range("D13").Value = "1.647.58"
$ws.Range("E13").Value = "  +1.49%  "

$ws.Range("D14").Value = "'4.17"
$ws.Range("E14").Value = "  +0.93%  "

$ws.Range("D15").Value = "'0.531"
$ws.Range("E15").Value = "  +1.61%  "

$ws.Range("D16").Value = "'65.04"
$ws.Range("E16").Value = "  +0.28%  "

$ws.Range("D17").Value = "26.706.19"

$ws.Range("D18").Value = "0.0₃0742"
$ws.Range("E18").Value = "  +0.36%  "

$ws.Range("D19").Value = "'217.08"
$ws.Range("E19").Value = "  +0.94%  "

$ws.Range("E20").Value = "  +0.23%  "

$ws.Range("D21").Value = "'4.35"
$ws.Range("E21").Value = "  +1.11%  "

$ws.Range("D22").Value = "'6.27"
$ws.Range("E22").Value = "  +0.09%  "

$ws.Range("B23").Value = "Toncoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D23").Value = "'2.29"
$ws.Range("E23").Value = "  +15.37%  "

$ws.Range("B24").Value = "Avalanche"
$ws.Range("C24").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D24").Value = "'9.49"
$ws.Range("E24").Value = "  +2.20%  "

$ws.Range("D25").Value = "'145.74"
$ws.Range("E25").Value = "  -1.56%  "

$ws.Range("E26").Value = "  +0.40%  "

$ws.Range("E27").Value = "  +0.06%  "

$ws.Range("D28").Value = "'7.11"
$ws.Range("E28").Value = "  +4.11%  "

$ws.Range("D29").Value = "'15.71"
$ws.Range("E29").Value = "  +0.90%  "

$ws.Range("D30").Value = "'0.0513"
$ws.Range("E30").Value = "  +0.89%  "

$ws.Range("D31").Value = "'1.17"
$ws.Range("E31").Value = "  +1.36%  "

$ws.Range("D32").Value = "'3.35"
$ws.Range("E32").Value = "  +1.10%  "

$ws.Range("D33").Value = "'3.00"
$ws.Range("E33").Value = "  +2.06%  "

$ws.Range("D34").Value = "1.274.61"
$ws.Range("E34").Value = "  +4.67%  "

$ws.Range("E35").Value = "  +3.49%  "

$ws.Range("E36").Value = "  +1.52%  "

$ws.Range("D37").Value = "'0.0178"
$ws.Range("E37").Value = "  +3.19%  "

$ws.Range("D38").Value = "'0.532"
$ws.Range("E38").Value = "  +5.60%  "

$ws.Range("D39").Value = "'0.820"
$ws.Range("E39").Value = "  +3.13%  "

$ws.Range("E40").Value = "  +0.26%  "

$ws.Range("D41").Value = "'0.814"
$ws.Range("E41").Value = "  +2.83%  "

$ws.Range("D42").Value = "'2.26"
$ws.Range("E42").Value = "  +0.29%  "

$ws.Range("D43").Value = "'5.44"

$ws.Range("D44").Value = "1.786.23"
$ws.Range("E44").Value = "  +1.24%  "

$ws.Range("D45").Value = "'91.53"
$ws.Range("E45").Value = "  -1.41%  "

$ws.Range("D46").Value = "'59.87"
$ws.Range("E46").Value = "  +9.16%  "

$ws.Range("E47").Value = "  +1.06%  "

$ws.Range("D48").Value = "'0.0515"
$ws.Range("E48").Value = "  +1.16%  "

$ws.Range("D49").Value = "'7.74"
$ws.Range("E49").Value = "  +3.38%  "

$ws.Range("D50").Value = "'0.0965"
$ws.Range("E50").Value = "  +1.46%  "

$ws.Range("E51").Value = "  +0.25%  "
